$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

# Row 2
Set-TextValue $ws.Range("D2") "30.771.64"
$ws.Range("E2").Value = "  -0.61%  "

# Row 3
Set-TextValue $ws.Range("D3") "1.936.50"
$ws.Range("E3").Value = "  -0.71%  "

# Row 4
$ws.Range("E4").Value = "  +0.10%  "

# Row 5
Set-TextValue $ws.Range("D5") "243.35"
$ws.Range("E5").Value = "  -0.72%  "

# Row 6
Set-TextValue $ws.Range("D6") "0.9998"
$ws.Range("E6").Value = "  +0.22%  "

# Row 7
Set-TextValue $ws.Range("D7") "0.4885"
$ws.Range("E7").Value = "  +0.15%  "

# Row 8
Set-TextValue $ws.Range("D8") "0.2952"
$ws.Range("E8").Value = "  -0.56%  "

# Row 9
Set-TextValue $ws.Range("D9") "0.06885"
$ws.Range("E9").Value = "  +0.83%  "

# Row 10
$ws.Range("E10").Value = "  +0.73%  "

# Row 11
Set-TextValue $ws.Range("D11") "104.98"
$ws.Range("E11").Value = "  -2.20%  "

# Row 12
$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextValue $ws.Range("D12") "1.941.66"
$ws.Range("E12").Value = "  -0.33%  "

# Row 13
$ws.Range("B13").Value = "TRON"
$ws.Range("C13").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
Set-TextValue $ws.Range("D13") "0.07782"
$ws.Range("E13").Value = "  +0.65%  "

# Row 14
Set-TextValue $ws.Range("D14") "5.349"
$ws.Range("E14").Value = "  -2.29%  "

# Row 15
Set-TextValue $ws.Range("D15") "0.7027"

# Row 16
Set-TextValue $ws.Range("D16") "272.78"
$ws.Range("E16").Value = "  -3.45%  "

# Row 17
Set-TextValue $ws.Range("D17") "30.794.88"
$ws.Range("E17").Value = "  -0.58%  "

# Row 18
$ws.Range("B18").Value = "Uniswap"
$ws.Range("C18").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
Set-TextValue $ws.Range("D18") "5.684"
$ws.Range("E18").Value = "  +3.34%  "

# Row 19
$ws.Range("B19").Value = "ShibaInu"
$ws.Range("C19").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
Set-TextValue $ws.Range("D19") "0.000007718"
$ws.Range("E19").Value = "  -0.08%  "

# Row 20
Set-TextValue $ws.Range("D20") "13.13"
$ws.Range("E20").Value = "  -0.88%  "

# Row 21
Set-TextValue $ws.Range("D21") "0.9992"
$ws.Range("E21").Value = "  +0.19%  "

# Row 22
Set-TextValue $ws.Range("D22") "0.9998"
$ws.Range("E22").Value = "  +0.07%  "

# Row 23
Set-TextValue $ws.Range("D23") "6.534"
$ws.Range("E23").Value = "  +0.60%  "

# Row 24
Set-TextValue $ws.Range("D24") "9.809"
$ws.Range("E24").Value = "  -0.31%  "

# Row 25
Set-TextValue $ws.Range("D25") "164.81"
$ws.Range("E25").Value = "  -2.60%  "

# Row 26
Set-TextValue $ws.Range("D26") "19.58"
$ws.Range("E26").Value = "  -1.98%  "

# Row 27
Set-TextValue $ws.Range("D27") "2.165"
$ws.Range("E27").Value = "  -2.27%  "

# Row 28
Set-TextValue $ws.Range("D28") "0.1037"
$ws.Range("E28").Value = "  -1.59%  "

# Row 29
Set-TextValue $ws.Range("D29") "1.384"
$ws.Range("E29").Value = "  -1.83%  "

# Row 30
Set-TextValue $ws.Range("D30") "4.675"
$ws.Range("E30").Value = "  +2.33%  "

# Row 31
Set-TextValue $ws.Range("D31") "1.559"
$ws.Range("E31").Value = "  -1.56%  "

# Row 32
Set-TextValue $ws.Range("D32") "4.424"
$ws.Range("E32").Value = "  -0.98%  "

# Row 33
Set-TextValue $ws.Range("D33") "0.04904"
$ws.Range("E33").Value = "  -0.97%  "

# Row 34
Set-TextValue $ws.Range("D34") "0.7602"
$ws.Range("E34").Value = "  -0.76%  "

# Row 35
Set-TextValue $ws.Range("D35") "1.151"
$ws.Range("E35").Value = "  -1.70%  "

# Row 36
Set-TextValue $ws.Range("D36") "0.9991"
$ws.Range("E36").Value = "  +0.17%  "

# Row 37
Set-TextValue $ws.Range("D37") "2.731"
$ws.Range("E37").Value = "  +0.19%  "

# Row 38
Set-TextValue $ws.Range("D38") "0.02009"
$ws.Range("E38").Value = "  -0.58%  "

# Row 39
Set-TextValue $ws.Range("D39") "79.60"
$ws.Range("E39").Value = "  +6.41%  "

# Row 40
Set-TextValue $ws.Range("D40") "2.668"
$ws.Range("E40").Value = "  -1.01%  "

# Row 41
Set-TextValue $ws.Range("D41") "6.496"
$ws.Range("E41").Value = "  -0.65%  "

# Row 42
Set-TextValue $ws.Range("D42") "2.085"
$ws.Range("E42").Value = "  -3.32%  "

# Row 43
Set-TextValue $ws.Range("D43") "0.9042"
$ws.Range("E43").Value = "  +2.33%  "

# Row 44
Set-TextValue $ws.Range("D44") "0.4449"
$ws.Range("E44").Value = "  -1.01%  "

# Row 45
Set-TextValue $ws.Range("D45") "108.52"
$ws.Range("E45").Value = "  -0.87%  "

# Row 46
Set-TextValue $ws.Range("D46") "7.873"
$ws.Range("E46").Value = "  -3.93%  "

# Row 47
Set-TextValue $ws.Range("D47") "0.9997"
$ws.Range("E47").Value = "  +0.13%  "

# Row 48
Set-TextValue $ws.Range("D48") "991.20"
$ws.Range("E48").Value = "  +1.03%  "

# Row 49
Set-TextValue $ws.Range("D49") "0.1251"
$ws.Range("E49").Value = "  -1.02%  "

# Row 50
Set-TextValue $ws.Range("D50") "36.30"
$ws.Range("E50").Value = "  +1.46%  "

# Row 51
Set-TextValue $ws.Range("D51") "9.256"
$ws.Range("E51").Value = "  -1.58%  "
